# Add files via upload
# Trazado de requisitos relleno al completo
#
# Fills the BDD traceability-matrix sheet with "x" marks across the
# requirement-coverage columns (C/D/E) for every row that now has a block
# assigned, and brings the BDD sheet to the front (it was "Pruebas" before).

$wb = $excel.ActiveWorkbook
$wsBDD = $wb.Worksheets.Item("BDD")
$wsPruebas = $wb.Worksheets.Item("Pruebas")

# Every cell below is currently blank (style-only, s="14"); mark it "x".
$cells = @(
    "C3", "D3", "E3",
    "C4", "D4",
    "C5", "D5", "E5",
    "C6", "D6",
    "C8", "D8", "E8",
    "C9", "D9",
    "C10", "D10", "E10",
    "C11", "D11", "E11",
    "C12", "D12", "E12",
    "C14", "D14", "E14",
    "C15", "D15", "E15",
    "C16", "E16",
    "C17", "D17",
    "C20", "E20",
    "C23", "D23", "E23",
    "D24", "E24",
    "C27", "D27", "E27",
    "C28", "D28",
    "C29", "D29",
    "D30", "E30",
    "C32", "D32",
    "C33", "D33",
    "C35", "D35", "E35",
    "C36", "D36", "E36",
    "C37", "D37",
    "C38", "E38",
    "C42", "E42",
    "C44", "D44",
    "C46", "D46", "E46",
    "D47", "E47",
    "C50", "D50", "E50",
    "C51", "D51", "E51",
    "D53", "E53",
    "C54", "D54",
    "C55", "D55", "E55",
    "D56", "E56",
    "C59", "D59", "E59",
    "C60", "D60",
    "C61", "D61",
    "D63", "E63",
    "C64", "D64", "E64",
    "C67", "D67", "E67",
    "C68", "D68", "E68",
    "C69", "D69", "E69",
    "C70", "D70", "E70",
    "D72", "E72"
)

foreach ($cell in $cells) {
    $wsBDD.Range($cell).Value = "x"
}

# Restore the saved selection on each sheet, then make BDD the active
# (front-most) sheet/tab, as it was the last one touched before saving.
[void]$wsPruebas.Range("E68").Select()
$wsBDD.Activate()
[void]$wsBDD.Range("C68").Select()
